# Added min oil pressure + rearrangement
$wb = $excel.ActiveWorkbook

# --- Typography sheet: fill in the Wildcard Characters column for the
#     "Default" typography (row 4) with the "0-9" wildcard set. ---
$wsTypography = $wb.Worksheets.Item("Typography")
$wsTypography.Range("G4").Value = "0-9"

# --- Translation sheet: the two "small" typography rows that describe the
#     oil-pressure readout (SingleUseId84 / SingleUseId85) are rearranged
#     to use the "Default" typography instead. ---
$wsTranslation = $wb.Worksheets.Item("Translation")
$wsTranslation.Range("C54").Value = "Default"
$wsTranslation.Range("C55").Value = "Default"

# --- New rows for the added "min oil pressure" text pair (value + unit
#     descriptor, following the same pattern as the other readouts). These
#     rows are brand new cells inside the Table8 ListObject range, so reset
#     their style back to Normal afterwards -- otherwise the table
#     auto-formats newly populated rows with a banding border that the
#     original file doesn't have. ---
$newRange = $wsTranslation.Range("B65:F66")

$wsTranslation.Range("B65").Value = "SingleUseId97"
$wsTranslation.Range("C65").Value = "Default"
$wsTranslation.Range("D65").Value = "Center"
$wsTranslation.Range("E65").Value = "<value>"
$wsTranslation.Range("F65").Value = "LTR"

$wsTranslation.Range("B66").Value = "SingleUseId98"
$wsTranslation.Range("C66").Value = "Default"
$wsTranslation.Range("D66").Value = "Left"
$wsTranslation.Range("E66").Value = "'0.0"
$wsTranslation.Range("F66").Value = "LTR"

$newRange.Style = "Normal"
